$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by one
# day (45177 -> 45178, i.e. 2023-09-08 -> 2023-09-09) for every data row
# (rows 2 through 135).
for ($r = 2; $r -le 135; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
